$p = $ppt.ActivePresentation

# Locate slide 19 ("EIP-7792: Verifiable logs") which holds the shapes we need to edit.
$s = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $cand = $p.Slides.Item($i)
    if ($cand.Shapes.Title.TextFrame.TextRange.Text -eq "EIP-7792: Verifiable logs") {
        $s = $cand
        break
    }
}

# Update "Block root" -> "Block number"
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        if ($shp.TextFrame.TextRange.Text -eq "Block root") {
            $shp.TextFrame.TextRange.Text = "Block number"
        }
        elseif ($shp.TextFrame.TextRange.Text -eq "Tx root") {
            $shp.TextFrame.TextRange.Text = "Tx index"
        }
    }
}

# Add new Cloud shape with EIP-7745 reference text.
$cloud = $s.Shapes.AddShape(9, 596.40, 48.66, 341.37, 79.29)
$cloud.Name = "Cloud 3"
$cloud.TextFrame.TextRange.Text = "EIP-7745: Two dimensional log filter data structure"
$cloud.TextFrame.TextRange.ParagraphFormat.Alignment = 2
